# Applies the "Updated symbol list on Wed Jan  4 15:18:28 UTC 2023 with GitHub
# Actions" edit: refreshed Price (D), Volume(1h) (E) and Hora (G) columns for
# the crypto table on Sheet1, rows 2-51.
#
# Values in this sheet are stored as *text* (e.g. "255.40", "4.10%", "15"),
# not numbers, so each one is written with a leading apostrophe (forces text
# entry, matching the original author's literal strings, including trailing
# zeros) and then the cell style is reset to "Normal" so no stray number format
# / quote-prefix formatting is left behind on a cell that had none before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: (worksheet row, worksheet column, new text value)
#   column 4 = D (Price), column 5 = E (Volume(1h)), column 7 = G (Hora)
$updates = @(
    @(2, 4, "255.40"),
    @(2, 5, "4.10%"),
    @(2, 7, "15"),
    @(3, 4, "28.13"),
    @(3, 5, "-4.28%"),
    @(3, 7, "15"),
    @(4, 4, "5.213"),
    @(4, 5, "-1.02%"),
    @(4, 7, "15"),
    @(5, 4, "0.05847"),
    @(5, 5, "2.13%"),
    @(5, 7, "15"),
    @(6, 4, "6.695"),
    @(6, 5, "0.96%"),
    @(6, 7, "15"),
    @(7, 4, "0.8702"),
    @(7, 5, "1.52%"),
    @(7, 7, "15"),
    @(8, 4, "0.9543"),
    @(8, 5, "11.86%"),
    @(8, 7, "15"),
    @(9, 4, "0.1408"),
    @(9, 5, "2.55%"),
    @(9, 7, "15"),
    @(10, 4, "0.07126"),
    @(10, 5, "0.58%"),
    @(10, 7, "15"),
    @(11, 4, "0.03210"),
    @(11, 5, "2.29%"),
    @(11, 7, "15"),
    @(12, 4, "0.09212"),
    @(12, 5, "-1.16%"),
    @(12, 7, "15"),
    @(13, 4, "0.001551"),
    @(13, 5, "1.30%"),
    @(13, 7, "15"),
    @(14, 4, "0.0006110"),
    @(14, 5, "1.92%"),
    @(14, 7, "15"),
    @(15, 4, "0.006016"),
    @(15, 5, "-1.16%"),
    @(15, 7, "15"),
    @(16, 4, "3.497"),
    @(16, 5, "-0.56%"),
    @(16, 7, "15"),
    @(17, 4, "3.215"),
    @(17, 5, "0.54%"),
    @(17, 7, "15"),
    @(18, 4, "2.226"),
    @(18, 5, "2.01%"),
    @(18, 7, "15"),
    @(19, 4, "0.3181"),
    @(19, 5, "0.68%"),
    @(19, 7, "15"),
    @(20, 4, "0.03448"),
    @(20, 5, "3.20%"),
    @(20, 7, "15"),
    @(21, 5, "0.42%"),
    @(21, 7, "15"),
    @(22, 4, "3.548"),
    @(22, 5, "1.66%"),
    @(22, 7, "15"),
    @(23, 4, "0.04187"),
    @(23, 5, "1.01%"),
    @(23, 7, "15"),
    @(24, 4, "0.1368"),
    @(24, 5, "-3.01%"),
    @(24, 7, "15"),
    @(25, 4, "0.001225"),
    @(25, 5, "0.32%"),
    @(25, 7, "15"),
    @(26, 4, "0.004542"),
    @(26, 5, "9.31%"),
    @(26, 7, "15"),
    @(27, 5, "-0.05%"),
    @(27, 7, "15"),
    @(28, 4, "0.0001469"),
    @(28, 5, "1.29%"),
    @(28, 7, "15"),
    @(29, 7, "15"),
    @(30, 7, "15"),
    @(31, 7, "15"),
    @(32, 7, "15"),
    @(33, 7, "15"),
    @(34, 7, "15"),
    @(35, 7, "15"),
    @(36, 7, "15"),
    @(37, 7, "15"),
    @(38, 7, "15"),
    @(39, 7, "15"),
    @(40, 4, "0.03821"),
    @(40, 5, "1.77%"),
    @(40, 7, "15"),
    @(41, 4, "0.005660"),
    @(41, 5, "-1.64%"),
    @(41, 7, "15"),
    @(42, 4, "0.1101"),
    @(42, 5, "3.28%"),
    @(42, 7, "15"),
    @(43, 4, "0.002356"),
    @(43, 5, "-2.56%"),
    @(43, 7, "15"),
    @(44, 4, "0.009740"),
    @(44, 5, "4.57%"),
    @(44, 7, "15"),
    @(45, 4, "0.00005405"),
    @(45, 5, "2.32%"),
    @(45, 7, "15"),
    @(46, 4, "0.00000000751"),
    @(46, 5, "0.11%"),
    @(46, 7, "15"),
    @(47, 4, "0.09015"),
    @(47, 5, "11.37%"),
    @(47, 7, "15"),
    @(48, 5, "-4.16%"),
    @(48, 7, "15"),
    @(49, 4, "0.00002104"),
    @(49, 5, "0.11%"),
    @(49, 7, "15"),
    @(50, 4, "0.0002004"),
    @(50, 5, "0.11%"),
    @(50, 7, "15"),
    @(51, 7, "15")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    # Leading apostrophe -> stored as literal text, preserving exact formatting
    $cell.Value = "'" + $val
    # Drop the quote-prefix style Excel adds for text that looks numeric, so the
    # cell keeps the same (default) style it had before the edit.
    $cell.Style = "Normal"
}
